$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 527, shifting existing rows 527:554 down to 528:555
$ws.Rows("527:527").Insert()

# Populate the newly inserted row 527 with the new record's data
$ws.Range("A527").Value = 10
$ws.Range("B527").Value = 'Vega Modelo de Temuco'
$ws.Range("C527").Value = 'La Araucanía'
$ws.Range("D527").Value = 45041
$ws.Range("E527").Value = 9
$ws.Range("F527").Value = 'Fruta'
$ws.Range("G527").Value = 100108
$ws.Range("H527").Value = 'Tropicales y subtropicales'
$ws.Range("I527").Value = 100108002
$ws.Range("J527").Value = 'Mango'
$ws.Range("K527").Value = 'Sin especificar'
$ws.Range("L527").Value = 'Primera'
$ws.Range("M527").Value = 185
$ws.Range("N527").Value = 8000
$ws.Range("O527").Value = 8000
$ws.Range("P527").Value = 8000
$ws.Range("Q527").Value = '$/bandeja 4 kilos'
$ws.Range("R527").Value = 'Perú'
$ws.Range("S527").Value = 2000
$ws.Range("T527").Value = 4
